# Scheduled runner update: refresh market-board derived profit figures
# across the FFXIV Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 68.8
$ws.Range("I6").Value = 68.8
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 206.4
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -94.39999999999998
$ws.Range("N6").ClearContents()

# Row 51
$ws.Range("H51").Value = 1688.1111
$ws.Range("I51").Value = 1400.5
$ws.Range("J51").Value = 1770.2858
$ws.Range("K51").Value = 1400.5
$ws.Range("L51").Value = 1770.2858
$ws.Range("M51").Value = -916.5
$ws.Range("N51").Value = -2738.2858

# Row 64
$ws.Range("H64").Value = 3260.6086
$ws.Range("J64").Value = 3574.875
$ws.Range("L64").Value = 3574.875
$ws.Range("N64").Value = -4070.875

# Row 67
$ws.Range("H67").Value = 3260.6086
$ws.Range("J67").Value = 3574.875
$ws.Range("L67").Value = 3574.875
$ws.Range("N67").Value = -5290.875

# Row 125
$ws.Range("H125").Value = 4190
$ws.Range("I125").Value = 1625
$ws.Range("K125").Value = 14625
$ws.Range("M125").Value = -12165

# Row 132
$ws.Range("H132").Value = 1666.6552
$ws.Range("I132").Value = 1651.2693
$ws.Range("J132").Value = 1800
$ws.Range("K132").Value = 4953.8079
$ws.Range("L132").Value = 5400
$ws.Range("M132").Value = -2423.8079
$ws.Range("N132").Value = -10460

# Row 137
$ws.Range("H137").Value = 8930898
$ws.Range("I137").Value = 17858806
$ws.Range("J137").Value = 2990.7144
$ws.Range("K137").Value = 53576418
$ws.Range("L137").Value = 8972.143199999999
$ws.Range("M137").Value = -53573868
$ws.Range("N137").Value = -14072.1432

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2428.1667
$ws.Range("I2").Value = 2541.6
$ws.Range("J2").Value = 1861
$ws.Range("K2").Value = 2541.6
$ws.Range("L2").Value = 1861
$ws.Range("M2").Value = -2428.6
$ws.Range("N2").Value = -2087

# Row 10
$ws.Range("H10").Value = 1380000
$ws.Range("I10").Value = 2750000
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 2750000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -2749830
$ws.Range("N10").Value = -10340

# Row 13
$ws.Range("H13").Value = 5000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5000
$ws.Range("N13").Value = -5288
$ws.Range("M13").ClearContents()

# Row 116
$ws.Range("H116").Value = 2428.1667
$ws.Range("I116").Value = 2541.6
$ws.Range("J116").Value = 1861
$ws.Range("K116").Value = 2541.6
$ws.Range("L116").Value = 1861
$ws.Range("M116").Value = -247.5999999999999
$ws.Range("N116").Value = -6449

# Row 132
$ws.Range("H132").Value = 2401.1843
$ws.Range("I132").Value = 2023.7812
$ws.Range("J132").Value = 4414
$ws.Range("K132").Value = 6071.3436
$ws.Range("L132").Value = 13242
$ws.Range("M132").Value = -3541.3436
$ws.Range("N132").Value = -18302

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2428.1667
$ws.Range("I3").Value = 2541.6
$ws.Range("J3").Value = 1861
$ws.Range("K3").Value = 2541.6
$ws.Range("L3").Value = 1861
$ws.Range("M3").Value = -2427.6
$ws.Range("N3").Value = -2089

# Row 15
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20454

# Row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# Row 94
$ws.Range("H94").Value = 696.1053000000001
$ws.Range("I94").Value = 585.5
$ws.Range("J94").Value = 885.7143
$ws.Range("K94").Value = 585.5
$ws.Range("L94").Value = 885.7143
$ws.Range("M94").Value = -134.5
$ws.Range("N94").Value = -1787.7143

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 18700
$ws.Range("J13").Value = 18700
$ws.Range("L13").Value = 18700
$ws.Range("N13").Value = -18978

# Row 99
$ws.Range("H99").Value = 2514.1333
$ws.Range("I99").Value = 1771.2
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 1771.2
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -273.2
$ws.Range("N99").Value = -6996

# Row 126
$ws.Range("H126").Value = 2514.1333
$ws.Range("I126").Value = 1771.2
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 5313.6
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -2843.6
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 2689.3333
$ws.Range("I36").Value = 1510
$ws.Range("J36").Value = 3118.182
$ws.Range("K36").Value = 4530
$ws.Range("L36").Value = 9354.545999999998
$ws.Range("M36").Value = -4361
$ws.Range("N36").Value = -9692.545999999998

# Row 68
$ws.Range("H68").Value = 982.4167
$ws.Range("I68").Value = 891.55554
$ws.Range("J68").Value = 1099.238
$ws.Range("K68").Value = 2674.66662
$ws.Range("L68").Value = 3297.714
$ws.Range("M68").Value = -1863.66662
$ws.Range("N68").Value = -4919.714

# Row 71
$ws.Range("H71").Value = 982.4167
$ws.Range("I71").Value = 891.55554
$ws.Range("J71").Value = 1099.238
$ws.Range("K71").Value = 8023.99986
$ws.Range("L71").Value = 9893.142
$ws.Range("M71").Value = -3967.99986
$ws.Range("N71").Value = -18005.142

# Row 107
$ws.Range("H107").Value = 45455316
$ws.Range("I107").Value = 202.6
$ws.Range("J107").Value = 83334580
$ws.Range("K107").Value = 607.8
$ws.Range("L107").Value = 250003740
$ws.Range("M107").Value = 1312.2
$ws.Range("N107").Value = -250007580

# Row 122
$ws.Range("H122").Value = 12386.277
$ws.Range("I122").Value = 23986.555
$ws.Range("J122").Value = 786
$ws.Range("K122").Value = 215878.995
$ws.Range("L122").Value = 7074
$ws.Range("M122").Value = -213428.995
$ws.Range("N122").Value = -11974

# Row 132
$ws.Range("H132").Value = 482171.38
$ws.Range("I132").Value = 813.9
$ws.Range("J132").Value = 919769.0600000001
$ws.Range("K132").Value = 7325.099999999999
$ws.Range("L132").Value = 8277921.540000001
$ws.Range("M132").Value = -4795.099999999999
$ws.Range("N132").Value = -8282981.540000001

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 913945.4
$ws.Range("I3").Value = 1433085.6
$ws.Range("J3").Value = 5450
$ws.Range("K3").Value = 1433085.6
$ws.Range("L3").Value = 5450
$ws.Range("M3").Value = -1432969.6
$ws.Range("N3").Value = -5682

# Row 13
$ws.Range("H13").Value = 2471
$ws.Range("I13").Value = 319.85715
$ws.Range("J13").Value = 10000
$ws.Range("K13").Value = 319.85715
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = -180.85715
$ws.Range("N13").Value = -10278

# Row 35
$ws.Range("H35").Value = 5013
$ws.Range("I35").Value = 5013
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 5013
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4715
$ws.Range("N35").ClearContents()

# Row 122
$ws.Range("H122").Value = 3699
$ws.Range("I122").Value = 4840
$ws.Range("J122").Value = 2985.875
$ws.Range("K122").Value = 14520
$ws.Range("L122").Value = 8957.625
$ws.Range("M122").Value = -12070
$ws.Range("N122").Value = -13857.625

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3190.6667
$ws.Range("I7").Value = 3112.5715
$ws.Range("J7").Value = 3300
$ws.Range("K7").Value = 3112.5715
$ws.Range("L7").Value = 3300
$ws.Range("M7").Value = -3000.5715
$ws.Range("N7").Value = -3524

# Row 61
$ws.Range("H61").Value = 13683.333
$ws.Range("I61").Value = 19025
$ws.Range("K61").Value = 19025
$ws.Range("M61").Value = -18823

# Row 113
$ws.Range("H113").Value = 13683.333
$ws.Range("I113").Value = 19025
$ws.Range("K113").Value = 19025
$ws.Range("M113").Value = -16855

# Row 126
$ws.Range("H126").Value = 3190.6667
$ws.Range("I126").Value = 3112.5715
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 9337.7145
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -6867.7145
$ws.Range("N126").Value = -14840

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 13666.667
$ws.Range("I54").Value = 13666.667
$ws.Range("K54").Value = 13666.667
$ws.Range("M54").Value = -13146.667

# Row 113
$ws.Range("H113").Value = 577.5454999999999
$ws.Range("I113").Value = 577.5454999999999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1732.6365
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 437.3635000000002
$ws.Range("N113").ClearContents()
